# Update the "interested count" (column F) values across the four worksheets
# to match the regenerated data snapshot (commit: "Update gh-pages to output
# generated at 456a3b4").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 549
$ws.Range("F3").Value = 983
$ws.Range("F4").Value = 74
$ws.Range("F7").Value = 1215
$ws.Range("F8").Value = 975
$ws.Range("F9").Value = 42
$ws.Range("F12").Value = 4347
$ws.Range("F13").Value = 584
$ws.Range("F14").Value = 151
$ws.Range("F15").Value = 1721
$ws.Range("F17").Value = 661
$ws.Range("F21").Value = 1100
$ws.Range("F22").Value = 1540
$ws.Range("F23").Value = 783
$ws.Range("F24").Value = 682
$ws.Range("F25").Value = 524
$ws.Range("F26").Value = 488
$ws.Range("F27").Value = 370
$ws.Range("F28").Value = 75
$ws.Range("F31").Value = 349
$ws.Range("F34").Value = 1458
$ws.Range("F36").Value = 11
$ws.Range("F38").Value = 4101

$ws = $wb.Worksheets.Item(2)
$ws.Range("F17").Value = 31
$ws.Range("F23").Value = 266
$ws.Range("F25").Value = 132
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = 20

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 1701
$ws.Range("F7").Value = 1054
$ws.Range("F8").Value = 148

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1701
$ws.Range("F5").Value = 1054
$ws.Range("F6").Value = 549
$ws.Range("F7").Value = 983
$ws.Range("F8").Value = 74
$ws.Range("F9").Value = 1215
$ws.Range("F10").Value = 975
$ws.Range("F12").Value = 42
$ws.Range("F17").Value = 148
$ws.Range("F19").Value = 4347
$ws.Range("F20").Value = 584
$ws.Range("F21").Value = 1721
$ws.Range("F22").Value = 661
$ws.Range("F27").Value = 1540
$ws.Range("F30").Value = 783
$ws.Range("F31").Value = 682
$ws.Range("F32").Value = 524
$ws.Range("F33").Value = 488
$ws.Range("F34").Value = 75
$ws.Range("F36").Value = 266
$ws.Range("F44").Value = 1458
$ws.Range("F45").Value = 11
$ws.Range("F49").Value = 4101
